# Sync attendance_reports: reorder "Recorded By" (column G) so that the
# literal token "System" (exact case) always appears first in the
# comma-separated list of recorders, preserving the relative order of the
# remaining entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $value = $cell.Value()

    if ($value -eq $null) {
        continue
    }

    $text = [string]$value
    if ($text -eq "") {
        continue
    }

    $parts = $text -split ","
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) {
            $hasSystem = $true
        }
    }

    if ($hasSystem) {
        $newParts = @("System")
        foreach ($p in $parts) {
            if (-not $p.Equals("System")) {
                $newParts += $p
            }
        }
        $newValue = $newParts -join ", "
        if (-not $newValue.Equals($text)) {
            $cell.Value = $newValue
        }
    }
}
